# Update "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-09-01 11:15:40"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-09-01 11:15:35"
$wsZhCn.Range("K2").Value = "2016-09-01 11:15:53"

# de-de sheet: "Correspond Handback DateTime" (K2)
$wsDeDe.Range("K2").Value = "2016-09-01 11:16:02"
